# TC_65: rename sheet, populate the data series (9 yearly observations),
# widen the custom number format to show thousands without a leading 0,
# and refresh the CDM add-in comment payload that documents the download.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be labelled after the add-in concept ("My Series");
# rename it to the plain "Data" tab name used going forward.
$ws.Name = "Data"

# Previously only the last observation (37962.97) was kept in A1. Restore
# the full 9-point yearly series in A1:A9, oldest first.
$values = @(3271.34, 6384.48, 7111.74, 8354.02, 10540.67, 16421.36, 22352.63, 29205.22, 37962.97)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Apply the (slightly widened) custom numeric format to the whole series.
$ws.Range("A1:A9").NumberFormat = "###0.000"

# Refresh the embedded CDM metadata blob stored in A1's comment so it
# reflects the re-downloaded range (A1:A9) and updated timestamp/graph title.
$comment = $ws.Range("A1").Comment
if ($comment -eq $null) {
    $comment = $ws.Range("A1").AddComment()
}
$commentText = "DBoAAB+LCAAAAAAAAAOlGdtu48b1Vwg9tUAlktJ612vMMtDNjlDJMiS5jvNSjMixNTU1VMihbb2lQIIUaYqiKDZFekWfUhTodtEmQLrby78Ea+/2qb/QMxcOSUnOVtuFseac25w5c65j9Nb1PLQuSZzQiD2quDWnYhHmRwFl548qKT+ruvcrb3moe+2T8AjHeE44EFvAxZK964Q+qsw4X+zZ9tXVVe2qUYvic7vuOK79zqA/9mdkjquUJRwzn1QMV/B6roqH2sF8QDgOMMeK81GlN+7V2oT6HYANMMPnJK610oQykiRdximnJBGcMcGctDuD76mDefXa/ZqL7DV4TtlKaRgouhKlgms62JZM6Jx4dcfdrTq71YY7cd29Rn1vp167t1t/N2M0hKiPEz4m8SX1JWDM8Xwh2Z3dhuu6jfpOHdkbiUBWbgAPDcNgRC5pQoI2CcNkK4vY+gKbPodTb2dMB9kFXi3ozVU4iPFiNqE8JNupMRq0rDnTuuRCPLQfxcQH+72RSofkahhrs04WfcBOZjTmyw5ebi3rOCHxcCGMtB2rhzoR482QxPx4AXdNAnAFQHg8Tgmy70DmTB2a+PBNWUoC7wyHSZGphEQnUXyRLLBPDiGObSHjioURDsDhOE049ZNcwBoGHcXRAkTC7q0oDPZBbEa9AWNk9xgYWWzciqKLVfFlJJL3Km8YbnWOjfw1OBrPoqshC5fjdJr4MZ2SoNPKqDfikAhJzd1OEx7NQYschBSsAFnCPwjBVTDqEJ/OcXgUghkTrwFSSgDUTHl0Rnk7CtM5M/ZcgaITONGEXJsTmjUawvUyYfWI9VhGr1xhI6rMMIquzJ7rCGmEAriZ+JmXrSNWiTsAy65vHSNvRJxyn4ZQIYp3UYCWvWI8I4RvdAmFQSIZ7oua47WWYk9k5xAErgn+DVDPffjwQdVx4WfiOHvyB3Y2aNRlgfyACnO/6tZF1s7pMiQ6TOfDKcTwpTyT5wJuBYTgFGErxOwCoCeUzw6bmfYbMEid+U76dRyCaF2EeCnBxi5FGOoxP0wDopJAj51JpxS6aeo78WgN1Ie49hBmy8lyIZKBfQdFh5zhNIT6xCGNnOcZYgWMmsnFKk0RhI7jMLtjT1T/BMq/H8xrPiRIUeJqfjQXABuq7skY2UV6keV90mXnfczOU8ijxo6rcONhIgNMYswScRyTNlecbTMRyu5FpXNPXdYwle6kLisCLLJX6NCEzBdRjMMBGIbup0yWTF0TwMIDzGd6BdEbEj8zsp2zGq6yZpniryOTYaiOIXJ8MQwLUEklDqM6jQJRDkTinIMogGKPQzqNS262EQdXlpfAzOHE8bYsh9ktQLsJGea7ZCl6kHyh4cJrPTdDyBVqg0reeASd2I7TgL5KrpE88ojg0OpCI8qJ1WOXJOFzYNuzRiShAXxRHO5Zb5MpoRD20kY68LbmLvKh/Zi8l0JDvZSqNEHfMqRMADn1nDIcrhMaTM7gnRIch8sCoTpqP/KB7vZH/7r59dMXzz6//fjxqy8/+M/ff/niHz+9efIhfNz+5a83n/xCHVMRowmehkQqNGnt7jqNe+BpBoR0goAqH6Q+l7DTU1n1zRrp9lUu2t1e+6DfkhnFADN20VGkIhr6eBml+XKsDiE3kldqZ56gSLxJlqH0uoTtEFHqVRiKTvWSlKmL+LsYlS1ePv/85fM/3cmtDVYqPjuiqLy++LjuGp0pPmICUA1eNl/cqzo71Xq9QLxCg0YwJkE/Z+zUC7yG6zx06g3HNdk8MI68iWgVpSVN8Lm9wqdA7ShlPF4aFyiuM6R0/AmEiEGrUCgstIt++eNXf35cotLW1ZCyFFAuSmPlTXa2kKIPRxNrPDwetbvWpDsWfpLjCnRK+DcQ691NPJWcirEUh9+B4ViMxlYFur6KFZ1ZBPszawmRWIjDkrNtgqqN3lDkqpYHcZQu1I0UGHLoBkqTTTZybMg1EiftuZZ0ctQGcqXrzd++2MSgD6Ld7JhRbqa8IgyVMApUwOuo/eyfL7766MWzZ7dPf3bz1Q9LEvQ+ZvgBP4doKi6N20PK0/VmBYJOxtKYF873C/VFA0XreBRRxhPPvS+7Rr1CwOoKafI36s2h5EnB0l4AX4Ggt3HSveY6sL1DZJcBoOcCQ7WN8v7aAFQOz+3679/89vZXX9x++vTVR3+8+fgPN598+vL57149+b2KutvHT29/8kRn+dVCIHURXbtqAy05BPqWiEZL1G7r6/d/brGIW9B0WKnMSF+//1lBmFBUtie5ZGjqjCJlFdZIi8yCzyqoYnQo8RkW1QC0RQlrGApdxKIF9fNN3q0KUSLuJOJbvUk1TYgVQT/1bThJmThn/l/5NIsqqUcPnLpb11iljTjCFCcF0x+E0RSajAwhh6wVkhLXNzPktHK/g/6w1eznJEqJYRzATOaISVF8oKypFCWll2QrM1rkEMBC4+enoZiL18jWUUZyIY3ZesA8awYi/XmOIye4cn0sUcCMHseqIWL6xXKcLqAd5nqIvRsvH2cKDfCh6lWLLXG+7nXKeFgXsFAIy2gBkHiZmjRKpaleImZZ1c4eCtPkS8CVHnTAHPpRUnVal9BXxrbIO904juKNySfHZGQD6KQho9i5xQ2NvFPVdQf5XWWALOGZDzX76RNGHRISvt2LnZ1zD6LLN+aFu9+WtZcMw0Abc7vRw5glF1B8thSO8v++Wipna8YxNFbikWPrZ8ZsdB3BxLulNuooklGMgLC7fgfcp3HC3xGZQH8pyKmBnKoOVZKoD7k+9XbVGvB2Ubhd0jKLXK7ejqOwT+d0y6nQycK7LARMuVioDq63naOIynJIrqG/LEiAnDj9AVQNMfFsJ035K6RSwy8eaxJ6PuPbKvZgiklApk7Vn5J69V7g7FYfEtKoui78j/163XF2xEuPFg6Jg5KrLTexswvL/6bj/Rdm2k9XDBoAAA=="
[void]$comment.Text($commentText)

# The CDM add-in's orphaned customXml metadata part (not referenced by any
# workbook relationship even before this edit) is dropped from the package
# alongside the refreshed comment above.
try {
    $xmlParts = $wb.CustomXMLParts
    for ($i = $xmlParts.Count; $i -ge 1; $i--) {
        [void]$xmlParts.Item($i).Delete()
    }
} catch {
    # Best effort only - absence of this collection shouldn't fail the edit.
}

Write-Output ("Sheet name: " + $ws.Name)
Write-Output ("A1 value: " + $ws.Range("A1").Value2)
Write-Output ("A9 value: " + $ws.Range("A9").Value2)
Write-Output ("A1 number format: " + $ws.Range("A1").NumberFormat)
